# Generate Report for Handback
# Updates the handback status / timestamps for file
# "1d810fc8-2299-46f8-a074-2c75f9035a63.md" (row 3) across the
# Overview, zh-cn, and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet (row 3 = 1d810fc8-...md) ---
$wsOverview.Range("E3").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("G3").Value = "2017-02-21 11:05:41"

# --- zh-cn sheet (row 3 = 1d810fc8-...md) ---
$wsZhCn.Range("C3").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("H3").Value = "2017-02-21 11:05:23"
$wsZhCn.Range("L3").Value = "2017-02-21 11:07:13"

# --- de-de sheet (row 3 = 1d810fc8-...md) ---
$wsDeDe.Range("C3").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("H3").Value = "2017-02-21 11:05:41"
$wsDeDe.Range("L3").Value = "2017-02-21 11:07:35"

$wb.Save()
